$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Linear
$ws.Range("B2").Value = 1.917154329115354
$ws.Range("C2").Value = 3.185411237504307
$ws.Range("D2").Value = -0.7667619768603215

# Row 3 - Decision Tree
$ws.Range("B3").Value = 1.264053467596442
$ws.Range("C3").Value = 1.312
$ws.Range("D3").Value = 0.6661037641116716

# Row 4 - Random Forest
$ws.Range("B4").Value = 1.282382578854724
$ws.Range("C4").Value = 1.344122424812032
$ws.Range("D4").Value = 0.6463121152129976

# Row 5 - Optimized Equation
$ws.Range("B5").Value = 1.500203407079482
$ws.Range("C5").Value = 1.948729800572808
$ws.Range("D5").Value = 0.337554223526647
